$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 46081 = 2026-02-28) for
# every data row (rows 2 through 391). Bump it by one day (serial 46082 =
# 2026-03-01) for all of them in one shot.
$lastRow = 391
$range = $ws.Range("C2:C$lastRow")
$range.Value = 46082
